# Generator update: add a per-speaker "hashcode" column and fix a typo
# ("Cleaveland" -> "Cleveland") in the Balto abstract, per commit:
#   "updated generator code to include hash"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Cleaveland" -> "Cleveland" typo in the Balto talk abstract (H22) ---
$baltoAbstract = @"
Almost 100 years ago, Balto, along with a pack of elite sled dogs helped save the community of Nome, Alaska, from a diphtheria outbreak. Balto is representative of the sled dogs of the era who were reputed for their hardiness, endurance, and tolerance of harsh conditions at a time when northern communities relied on them for transportation, protection, and companionship. Unlike modern breeds, Balto and his contemporaries were products of generations of mating of diverse, outbred dogs selected for their performance and hardiness.  Today, Balto is immortalized with a statue in Central Park and is physically preserved and on display at the Cleveland Museum of Natural History. We asked ourselves, could we use advances in genomic technologies and ancient DNA recovery techniques to parse through Balto’s DNA and understand selection emphasis and diversity occurring in sled dogs in the 1930’s? Could we decipher Balto’s physical appearance from his DNA and identify healthy developmental gene variants that may have endowed sled dogs in that era with abilities to withstand and thrive in harsh environments? Could we use our findings to guide modern breeders to breed healthier dogs? <br><br>Also, here's a photo of Dr. Huson racing as a musher in Alaska!<div class="box"> <img src="images/huson_sled.jpg" style = 'width: 550' class='center'> </div>
"@
$ws.Range("H22").Value = $baltoAbstract

# --- Add the new "hashcode" column (I) with a short slug per speaker ---
$hashcodes = @(
    "hashcode",
    "sethkorproski",
    "zachulibarri",
    "evesnyder",
    "jamesnagy",
    "benfried",
    "ligiacoelho",
    "meganbarrington",
    "sethstrickland",
    "danielasamur",
    "lukekeller",
    "hunteradams",
    "michaelcaporizzo",
    "mollyryan",
    "karlsmolenski",
    "daisyrosas",
    "alisonritterhaus",
    "zachulibarri2",
    "adamhawkins",
    "jayleeming",
    "danielsprocket",
    "heatherhuson",
    "alitahoward"
)

for ($i = 0; $i -lt $hashcodes.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 9).Value = $hashcodes[$i]
}

# --- Widen the speaker/position/flavor/title columns to fit the new content ---
$ws.Columns.Item(4).ColumnWidth = 13.1666666666667   # D: speaker      -> width 14
$ws.Columns.Item(5).ColumnWidth = 16.1666666666667   # E: position     -> width 17
$ws.Columns.Item(6).ColumnWidth = 16.0221354166667   # F: flavor       -> width ~16.86
$ws.Columns.Item(7).ColumnWidth = 26.0221354166667   # G: title        -> width ~26.86

# --- Leave the selection where the editor ended up after the edit ---
$ws.Range("H28").Select()

Write-Output "edit complete"
